# Update bug tracker statuses from "New" to "Resolved"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column E holds "Status" for each bug row (rows 2-10); mark all rows as Resolved
$ws.Range("E2:E10").Value = "Resolved"

# Update the active cell selection to E2, matching the saved view state
$ws.Range("E2").Select()

